$wb = $excel.ActiveWorkbook

# --- Sheet "展览": update F column (想去人数) values ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 563
$ws1.Range("F3").Value = 183
$ws1.Range("F4").Value = 329
$ws1.Range("F5").Value = 402
$ws1.Range("F6").Value = 256
$ws1.Range("F7").Value = 2372
$ws1.Range("F8").Value = 397
$ws1.Range("F9").Value = 6055
$ws1.Range("F10").Value = 154
$ws1.Range("F11").Value = 389
$ws1.Range("F12").Value = 18

# --- Sheet "演出": shift rows 3,4,5 content up into rows 2,3,4; then delete trailing row 5 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("B2").Value = '2024-07-18'
$ws2.Range("C2").Value = '南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《胡桃夹子》'
$ws2.Range("D2").Value = '龙堤路25号 广西文化艺术中心'
$ws2.Range("E2").Value = '2024.07.18 20:00-07.18 21:30'
$ws2.Range("F2").Value = 11
$ws2.Range("G2").Value = 108
$ws2.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=85816'
$ws2.Range("I2").Value = '//i0.hdslb.com/bfs/openplatform/202405/SN0ZyGVj1715675672714.jpeg'
$ws2.Range("B3").Value = '2024-07-19'
$ws2.Range("C3").Value = '南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《天鹅湖》 '
$ws2.Range("D3").Value = '龙堤路25号 广西文化艺术中心'
$ws2.Range("E3").Value = '2024.07.19 20:00-07.19 22:00'
$ws2.Range("F3").Value = 15
$ws2.Range("G3").Value = 108
$ws2.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=85831'
$ws2.Range("I3").Value = '//i1.hdslb.com/bfs/openplatform/202405/ZyyeeOUo1715677877362.jpeg'
$ws2.Range("B4").Value = '2024-08-10'
$ws2.Range("C4").Value = '南宁·限时7折|浪漫七夕《一生所爱》《爱乐之城》《假如爱有天意》经典浪漫电影主题音乐会'
$ws2.Range("D4").Value = '龙堤路25号 广西文化艺术中心'
$ws2.Range("E4").Value = '2024.08.10 20:00-08.10 21:30'
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = 99
$ws2.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=87729'
$ws2.Range("I4").Value = '//i1.hdslb.com/bfs/openplatform/202406/qKUDMYOh1718177639735.png'
$ws2.Rows(5).Delete()

# --- Sheet "全部类型": shift rows 3..16 content up into rows 2..15; then delete trailing row 16 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("B2").Value = '2024-07-06'
$ws4.Range("C2").Value = '南宁·小蜜蜂动漫嘉年华2.0'
$ws4.Range("D2").Value = '亭洪路45号 百益上河城'
$ws4.Range("E2").Value = '2024.07.06 10:00-07.06 17:00'
$ws4.Range("F2").Value = 563
$ws4.Range("G2").Value = 50
$ws4.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=84925'
$ws4.Range("I2").Value = '//i2.hdslb.com/bfs/openplatform/202404/YjFyyYq51713508727131.jpeg'
$ws4.Range("B3").Value = '2024-07-06'
$ws4.Range("C3").Value = '南宁·首届童话梦境Lolita茶会'
$ws4.Range("D3").Value = '明秀东路157号 利泰国际大酒店'
$ws4.Range("E3").Value = '2024.07.06 13:00-07.06 17:00'
$ws4.Range("F3").Value = 183
$ws4.Range("G3").Value = 88
$ws4.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=85776'
$ws4.Range("I3").Value = '//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg'
$ws4.Range("B4").Value = '2024-07-12'
$ws4.Range("C4").Value = '南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展'
$ws4.Range("D4").Value = '民族大道106号 南宁国际会展中心'
$ws4.Range("E4").Value = '2024.07.12 09:30-07.14 17:00'
$ws4.Range("F4").Value = 329
$ws4.Range("G4").Value = 50
$ws4.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=87182'
$ws4.Range("I4").Value = '//i1.hdslb.com/bfs/openplatform/202406/x4UZPn301718159475475.jpeg'
$ws4.Range("B5").Value = '2024-07-13'
$ws4.Range("C5").Value = '南宁·0713国乙ONLY'
$ws4.Range("D5").Value = '亭洪路45号 水明漾宴会中心'
$ws4.Range("E5").Value = '2024.07.13 09:30-07.13 21:00'
$ws4.Range("F5").Value = 402
$ws4.Range("G5").Value = 68
$ws4.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=86378'
$ws4.Range("I5").Value = '//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg'
$ws4.Range("B6").Value = '2024-07-14'
$ws4.Range("C6").Value = '广西·首届明日方舟only展 - 花庭圣梦'
$ws4.Range("D6").Value = '明秀东路157号 利泰国际大酒店'
$ws4.Range("E6").Value = '2024.07.14 09:00-07.14 18:00'
$ws4.Range("F6").Value = 256
$ws4.Range("G6").Value = 69
$ws4.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=85852'
$ws4.Range("I6").Value = '//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg'
$ws4.Range("B7").Value = '2024-07-18'
$ws4.Range("C7").Value = '南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《胡桃夹子》'
$ws4.Range("D7").Value = '龙堤路25号 广西文化艺术中心'
$ws4.Range("E7").Value = '2024.07.18 20:00-07.18 21:30'
$ws4.Range("F7").Value = 11
$ws4.Range("G7").Value = 108
$ws4.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=85816'
$ws4.Range("I7").Value = '//i0.hdslb.com/bfs/openplatform/202405/SN0ZyGVj1715675672714.jpeg'
$ws4.Range("B8").Value = '2024-07-19'
$ws4.Range("C8").Value = '南宁·限时6折|俄罗斯圣彼得堡古典芭蕾舞剧院《天鹅湖》 '
$ws4.Range("D8").Value = '龙堤路25号 广西文化艺术中心'
$ws4.Range("E8").Value = '2024.07.19 20:00-07.19 22:00'
$ws4.Range("F8").Value = 15
$ws4.Range("G8").Value = 108
$ws4.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=85831'
$ws4.Range("I8").Value = '//i1.hdslb.com/bfs/openplatform/202405/ZyyeeOUo1715677877362.jpeg'
$ws4.Range("B9").Value = '2024-07-20'
$ws4.Range("C9").Value = '南宁·AB动漫游戏嘉年华'
$ws4.Range("D9").Value = '三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心'
$ws4.Range("E9").Value = '2024.07.20 09:30-07.21 17:00'
$ws4.Range("F9").Value = 2372
$ws4.Range("G9").Value = 60
$ws4.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=84862'
$ws4.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg'
$ws4.Range("B10").Value = '2024-07-20'
$ws4.Range("C10").Value = '横州·第二届海棠动漫游戏嘉年华'
$ws4.Range("D10").Value = '茉莉花大道 横州国际大酒店'
$ws4.Range("E10").Value = '2024.07.20 09:30-07.20 17:00'
$ws4.Range("F10").Value = 397
$ws4.Range("G10").Value = 30
$ws4.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=84799'
$ws4.Range("I10").Value = '//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg'
$ws4.Range("B11").Value = '2024-07-27'
$ws4.Range("C11").Value = '南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）'
$ws4.Range("D11").Value = '民族大道106号 南宁国际会展中心'
$ws4.Range("E11").Value = '2024.07.27 09:30-07.28 17:30'
$ws4.Range("F11").Value = 6055
$ws4.Range("G11").Value = 55
$ws4.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=85264'
$ws4.Range("I11").Value = '//i1.hdslb.com/bfs/openplatform/202406/JxFed5iv1718622152091.jpeg'
$ws4.Range("B12").Value = '2024-08-03'
$ws4.Range("C12").Value = '南宁·火影忍者only'
$ws4.Range("D12").Value = '厢竹大道65号 桔子酒店'
$ws4.Range("E12").Value = '2024.08.03 10:00-08.03 17:00'
$ws4.Range("F12").Value = 154
$ws4.Range("G12").Value = 68
$ws4.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=86994'
$ws4.Range("I12").Value = '//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg'
$ws4.Range("B13").Value = '2024-08-03'
$ws4.Range("C13").Value = '南宁·蔚蓝档案only'
$ws4.Range("D13").Value = '亭洪路45号 百益上河城'
$ws4.Range("E13").Value = '2024.08.03 09:00-08.03 17:00'
$ws4.Range("F13").Value = 389
$ws4.Range("G13").Value = 68
$ws4.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=85370'
$ws4.Range("I13").Value = '//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg'
$ws4.Range("B14").Value = '2024-08-10'
$ws4.Range("C14").Value = '南宁·限时7折|浪漫七夕《一生所爱》《爱乐之城》《假如爱有天意》经典浪漫电影主题音乐会'
$ws4.Range("D14").Value = '龙堤路25号 广西文化艺术中心'
$ws4.Range("E14").Value = '2024.08.10 20:00-08.10 21:30'
$ws4.Range("F14").Value = 0
$ws4.Range("G14").Value = 99
$ws4.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=87729'
$ws4.Range("I14").Value = '//i1.hdslb.com/bfs/openplatform/202406/qKUDMYOh1718177639735.png'
$ws4.Range("B15").Value = '2024-11-02'
$ws4.Range("C15").Value = '南宁·万圣漫控嘉年华10'
$ws4.Range("D15").Value = '亭洪路45号 百益上河城'
$ws4.Range("E15").Value = '2024.11.02 11:00-11.03 22:00'
$ws4.Range("F15").Value = 18
$ws4.Range("G15").Value = 50
$ws4.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=87820'
$ws4.Range("I15").Value = '//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg'
$ws4.Rows(16).Delete()
